$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "High"
$ws.Range("F3").Value = "Low"
$ws.Range("F4").Value = "Medium"
$ws.Range("F5").Value = "High"

$ws.Range("A1:F5").Select()
